# "project is updated till add new customer"
# The TestCaseID column (A) is removed from the LoginDetails sheet; the
# Username column (previously B) becomes A and the Password column
# (previously C) becomes B. The TestCaseID/TC_00x shared strings become
# unused and drop out of the saved workbook automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift Username (col B) into col A and Password (col C) into col B,
# row by row, for the 6 rows of data (header + 5 records).
for ($r = 1; $r -le 6; $r++) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 3).Value2
}

# Remove the now-duplicated old column C content (old Password column).
$ws.Range("C1:C6").Clear()

# Leave the sheet with the new A1:B6 table selected, like the source file.
$ws.Range("A1:B6").Select()
